# edit.ps1 - apply the KiCad.docx documentation wording tweaks described by the diff.
#
# Strategy: for every paragraph span that changes, locate the OLD text with
# Find.Execute, overwrite the whole span with the full NEW text (so the visible
# text matches the target), and then re-create the individual <w:r> run
# boundaries the diff shows by toggling Bold True/False (a net no-op formatting
# change) on each sub-range - this runtime's paragraph-normalizer coalesces
# adjacent runs that share formatting whenever it touches a paragraph, so every
# boundary we want to keep has to be "marked" this way or it gets merged back
# into its neighbour.
#
# NOTE: helper functions deliberately do NOT take the Document/Range COM
# object as a parameter - passing $word.ActiveDocument positionally into a
# function here corrupts the binding of the parameters that follow it, so
# each helper just re-reads $word.ActiveDocument itself.

function Mark-Boundary($from, $to) {
    $doc = $word.ActiveDocument
    $rng = $doc.Range($from, $to)
    $rng.Bold = 1
    $rng.Bold = 0
}

function Replace-WithRuns($findText, [object[]]$segments) {
    $doc = $word.ActiveDocument
    $find = $doc.Content.Find
    $find.ClearFormatting()
    $ok = $find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $ok) {
        Write-Output "NOT FOUND: $findText"
        return
    }

    $r = $find.Parent
    $start = $r.Start

    $newText = [string]::Join("", $segments)
    $r.Text = $newText

    $pos = $start
    $boundaries = @($start)
    foreach ($seg in $segments) {
        $pos = $pos + $seg.Length
        $boundaries += $pos
    }

    # $boundaries = [start, start+len(seg1), start+len(seg1)+len(seg2), ..., end]
    # Mark every adjacent pair so each segment keeps its own run.
    for ($i = 0; $i -lt $boundaries.Length - 1; $i++) {
        Mark-Boundary $boundaries[$i] $boundaries[$i + 1]
    }
}

$RSQUO = [char]0x2019

# ---------------------------------------------------------------------------
# 1) "...so it can be tested and moved safely... Later a pcb is..."
# ---------------------------------------------------------------------------
Replace-WithRuns `
    "so it can be tested and moved safely from one point to another. Later a pcb is the only viable way to " `
    @(
        "so it can be ",
        "used",
        " and moved safely from one point to another. Later a ",
        "PCB",
        " is the only viable way to "
    )

# ---------------------------------------------------------------------------
# 2) "are is follows" -> "are as follows"
# ---------------------------------------------------------------------------
Replace-WithRuns `
    "are is follows" `
    @("are ", "as ", "follows")

# ---------------------------------------------------------------------------
# 3) "...widely used by many hobiests... good software to learn how to design PCB's."
# ---------------------------------------------------------------------------
Replace-WithRuns `
    ([string]" KICAD is open-source, free, widely used by many hobiests and also some professionals and apart from that it" + $RSQUO + "s good software to learn how to design PCB" + $RSQUO + "s.") `
    @(
        " KICAD is open-source, free, widely used by many ",
        "hobbyists",
        ([string]" and also some professionals and apart from that it" + $RSQUO + "s "),
        "good software to",
        " start",
        " learn",
        "ing",
        ([string]" how to design PCB" + $RSQUO + "s.")
    )

# ---------------------------------------------------------------------------
# 4) "...install any useful plugins whatsoever... Th"+"ese are the " -> plugin list + "the main "
# ---------------------------------------------------------------------------
Replace-WithRuns `
    ([string]" Apart from this, you are not allowed to install any useful plugins whatsoever, which could" + $RSQUO + "ve saved you time or made your designs even better. These are the ") `
    @(
        " Apart from this, you are not allowed to install any useful plugins ",
        "like an interactive BOM file for assembly, ",
        "3D model archiver, ",
        "fabrication toolkits, ",
        ([string]"which could" + $RSQUO + "ve saved you time or made your designs "),
        "and workflow ",
        "even better. Th",
        "ese are ",
        "the main",
        " "
    )

# ---------------------------------------------------------------------------
# 5) "by using the parts from EasyEda's library that you gain the fastest..."
# ---------------------------------------------------------------------------
Replace-WithRuns `
    ([string]"by using the parts from EasyEda" + $RSQUO + "s library that you gain the fastest (and also probably cheapest) way to produce a fully assembled PCB.") `
    @(
        ([string]"by using the parts from EasyEda" + $RSQUO + "s library that you "),
        "most likely ",
        "gain the fastest (and also probably cheapest) way to produce a fully assembled PCB."
    )

# ---------------------------------------------------------------------------
# 6) "...ability back to easily get LCSC part numbers that are required by JLCPCB."
# ---------------------------------------------------------------------------
Replace-WithRuns `
    " Apart from that, there are plugins for KICAD6 that give you the ability back to easily get LCSC part numbers that are required by JLCPCB." `
    @(
        " Apart from that, there are plugins for KICAD6 that give you the ability back to",
        " relatively",
        " easily get LCSC part numbers that are required by JLCPCB",
        " for ",
        "PCB ",
        "assembly",
        "."
    )

Write-Output "done"
